$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mark cm014's two entries (rows 14 and 15) as linked/covered ("link_it" column C)
$ws.Range("C14").Value = $true
$ws.Range("C15").Value = $true

# Update the active selection to reflect where the user was working (C16)
$ws.Range("C16").Select()
